$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.0884
$ws.Range("C6").Value = -12.9497
$ws.Range("D10").Value = -7.963499999999995
$ws.Range("A14").Value = -21.91079999999999
$ws.Range("B15").Value = 4.917999999999997
$ws.Range("D15").Value = -8.5518
$ws.Range("A16").Value = -22.27420000000001
$ws.Range("C18").Value = -14.25759999999999
$ws.Range("D18").Value = -9.286200000000001
$ws.Range("C19").Value = -12.6077
$ws.Range("A21").Value = -21.83559999999998
$ws.Range("B21").Value = 5.305599999999998
$ws.Range("D21").Value = -7.821699999999996
$ws.Range("B22").Value = 8.557400000000007
$ws.Range("D22").Value = -8.432600000000003
$ws.Range("A23").Value = -20.04219999999998
$ws.Range("B24").Value = 5.626299999999999
$ws.Range("D24").Value = -6.60439999999999
$ws.Range("A25").Value = -21.8004
$ws.Range("A26").Value = -21.15619999999998
$ws.Range("B27").Value = 6.657900000000007
$ws.Range("B28").Value = 6.269600000000002
$ws.Range("A29").Value = -21.06479999999997
$ws.Range("D33").Value = -8.426499999999999
$ws.Range("C35").Value = -13.49270000000001
$ws.Range("B36").Value = 8.773399999999999
$ws.Range("B39").Value = 8.747299999999999
$ws.Range("A40").Value = -20.16740000000001
$ws.Range("C44").Value = -13.2634
$ws.Range("B45").Value = 4.652800000000004
$ws.Range("D46").Value = -7.875799999999997
$ws.Range("C47").Value = -12.1928
$ws.Range("B48").Value = 5.382900000000004
$ws.Range("B49").Value = 5.2812
$ws.Range("D49").Value = -8.255000000000001
$ws.Range("C50").Value = -13.72909999999998
$ws.Range("C51").Value = -11.4125
$ws.Range("B52").Value = 5.164300000000002
$ws.Range("C52").Value = -10.92899999999999
$ws.Range("A53").Value = -20.10539999999999
$ws.Range("B53").Value = 10.012
$ws.Range("B54").Value = 4.791700000000005
$ws.Range("C55").Value = -14.1878
$ws.Range("D56").Value = -8.344899999999996
$ws.Range("A57").Value = -20.31779999999997
$ws.Range("B57").Value = 8.442100000000003
$ws.Range("C57").Value = -12.4557
$ws.Range("C58").Value = -13.2719
$ws.Range("A59").Value = -22.96990000000002
$ws.Range("D61").Value = -8.345899999999997
$ws.Range("C64").Value = -10.4686
$ws.Range("A65").Value = -21.7937
$ws.Range("D65").Value = -8.306999999999995
$ws.Range("C66").Value = -12.883
$ws.Range("D66").Value = -7.763600000000001
$ws.Range("A69").Value = -21.6225
$ws.Range("B70").Value = 4.820800000000002
$ws.Range("B71").Value = 4.603699999999995
$ws.Range("D74").Value = -8.520000000000005
$ws.Range("D75").Value = -8.1614
$ws.Range("D77").Value = -6.534599999999998
$ws.Range("A79").Value = -20.2147
$ws.Range("C80").Value = -12.8693
$ws.Range("A83").Value = -21.83659999999999
$ws.Range("C83").Value = -12.4644
$ws.Range("B86").Value = 5.411500000000002
$ws.Range("B87").Value = 5.566199999999998
$ws.Range("D87").Value = -9.0731
$ws.Range("D88").Value = -7.794899999999999
$ws.Range("B89").Value = 5.108899999999998
$ws.Range("A91").Value = -20.37789999999998
$ws.Range("C92").Value = -10.2098
$ws.Range("A93").Value = -21.22620000000001
$ws.Range("C94").Value = -10.63899999999999
$ws.Range("C96").Value = -10.34440000000001
$ws.Range("C97").Value = -10.9583
$ws.Range("A100").Value = -22.0933
$ws.Range("D100").Value = -7.977899999999997
$ws.Range("B101").Value = 6.793800000000007
$ws.Range("C101").Value = -12.2701
$ws.Range("D101").Value = -8.432699999999993
$ws.Range("A103").Value = -22.11950000000001
